$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 86.29678392075563
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 110.1084396136492

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 6.741336633845642
